$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PEIIR-EFPEIIR")

# The single "natural gas nonpeaker" fuel column is being split into two
# separate fuel columns: "natural gas steam turbine" and
# "natural gas combined cycle". Insert a new column after the existing
# "natural gas nonpeaker" column (column C) to make room for the second
# new fuel column.
$ws.Columns.Item(4).Insert()

# Rename the (former "natural gas nonpeaker") column C header and give the
# newly inserted column D its header.
$ws.Cells.Item(1, 3).Value = "natural gas steam turbine"
$ws.Cells.Item(1, 4).Value = "natural gas combined cycle"

# Fill the new column's data rows (2-13) with 0, matching every other
# fuel's improvement-rate value on this sheet.
for ($r = 2; $r -le 13; $r++) {
    $ws.Cells.Item($r, 4).Value = 0
}

# Match formatting: new column gets the same width as the column it split
# from, header row becomes shorter (30 instead of 75) now that there's an
# extra column carrying part of the former header text, and column A is
# widened now that the row is short.
$ws.Columns.Item(4).ColumnWidth = $ws.Columns.Item(3).ColumnWidth
$ws.Rows.Item(1).RowHeight = 30
$ws.Columns.Item(1).ColumnWidth = 17
